$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for cells A2:A201 (in row order)
$values = @(3,2,3,3,3,3,1,3,1,3,2,1,1,1,1,2,1,3,1,3,2,3,1,2,1,3,1,3,3,1,3,3,1,3,3,2,3,3,1,1,3,3,3,3,3,1,3,1,3,3,3,1,1,1,2,1,2,3,1,3,3,2,2,1,2,1,1,1,3,2,1,3,3,1,2,3,3,3,3,1,1,3,3,3,3,1,3,1,2,3,2,2,3,3,3,3,3,3,3,2,1,1,1,1,3,3,3,1,2,3,2,3,3,3,3,3,3,3,2,3,3,3,3,3,2,1,1,3,3,3,2,3,2,2,3,1,1,3,3,1,2,3,1,3,2,2,1,3,2,3,2,1,3,1,2,1,3,3,3,1,1,3,2,3,1,1,3,3,2,3,2,1,3,3,1,3,2,2,3,3,2,1,3,3,3,3,3,3,3,1,3,3,3,3,3,3,2,2,2,3)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Remove the now-extra rows (202 through 251) that existed in the original data
$lastOldRow = 251
$firstRemoveRow = 202
if ($lastOldRow -ge $firstRemoveRow) {
    $deleteRange = $ws.Range("A" + $firstRemoveRow + ":A" + $lastOldRow).EntireRow
    $deleteRange.Delete()
}
